$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each written value to remain plain text (matching the source
# workbooks inlineStr cell type) instead of being auto-detected as a
# number by Excels input parser, while leaving cell styling untouched.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "43.724.21"
Set-TextValue "E2" "  -0.16%  "
Set-TextValue "D3" "2.285.96"
Set-TextValue "E3" "  +3.29%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "B5" "Solana"
Set-TextValue "C5" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D5" "94.79"
Set-TextValue "E5" "  +8.63%  "
Set-TextValue "B6" "BNB"
Set-TextValue "C6" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D6" "267.79"
Set-TextValue "E6" "  +2.47%  "
Set-TextValue "D7" "0.622"
Set-TextValue "E7" "  +0.98%  "
Set-TextValue "E8" "  -0.05%  "
Set-TextValue "D9" "0.623"
Set-TextValue "E9" "  +3.63%  "
Set-TextValue "D10" "46.29"
Set-TextValue "E10" "  +2.70%  "
Set-TextValue "D11" "0.0932"
Set-TextValue "E11" "  +1.48%  "
Set-TextValue "D12" "8.07"
Set-TextValue "E12" "  +7.62%  "
Set-TextValue "E13" "  +0.43%  "
Set-TextValue "D14" "2.624.31"
Set-TextValue "E14" "  +2.98%  "
Set-TextValue "D15" "15.35"
Set-TextValue "E15" "  +5.65%  "
Set-TextValue "D16" "0.836"
Set-TextValue "E16" "  +7.14%  "
Set-TextValue "D17" "2.264.05"
Set-TextValue "E17" "  +2.35%  "
Set-TextValue "D18" "43.707.60"
Set-TextValue "E18" "  -0.09%  "
Set-TextValue "E19" "  +1.57%  "
Set-TextValue "D20" "6.21"
Set-TextValue "E20" "  +4.39%  "
Set-TextValue "D21" "70.95"
Set-TextValue "E21" "  +1.55%  "
Set-TextValue "D22" "2.29"
Set-TextValue "E22" "  -1.90%  "
Set-TextValue "D23" "10.02"
Set-TextValue "E23" "  +11.91%  "
Set-TextValue "D24" "235.17"
Set-TextValue "E24" "  +1.45%  "
Set-TextValue "E25" "  +0.09%  "
Set-TextValue "D26" "11.33"
Set-TextValue "E26" "  +6.52%  "
Set-TextValue "D27" "2.47"
Set-TextValue "E27" "  +9.65%  "
Set-TextValue "D28" "3.42"
Set-TextValue "E28" "  -3.55%  "
Set-TextValue "D29" "39.72"
Set-TextValue "E29" "  -0.15%  "
Set-TextValue "E30" "  -0.40%  "
Set-TextValue "D31" "22.12"
Set-TextValue "E31" "  +8.04%  "
Set-TextValue "D32" "172.13"
Set-TextValue "E32" "  -1.46%  "
Set-TextValue "D33" "0.0913"
Set-TextValue "E33" "  +4.33%  "
Set-TextValue "D34" "5.57"
Set-TextValue "E34" "  +3.16%  "
Set-TextValue "D35" "0.125"
Set-TextValue "E35" "  +1.56%  "
Set-TextValue "E36" "  +0.46%  "
Set-TextValue "D37" "4.46"
Set-TextValue "E37" "  -0.56%  "
Set-TextValue "D38" "0.0348"
Set-TextValue "E38" "  -2.77%  "
Set-TextValue "D39" "3.42"
Set-TextValue "E39" "  +13.66%  "
Set-TextValue "D40" "0.241"
Set-TextValue "E40" "  +20.45%  "
Set-TextValue "D41" "2.27"
Set-TextValue "E41" "  +8.68%  "
Set-TextValue "D42" "12.38"
Set-TextValue "E42" "  -0.86%  "
Set-TextValue "D43" "1.32"
Set-TextValue "E43" "  +17.57%  "
Set-TextValue "D44" "5.44"
Set-TextValue "E44" "  -1.64%  "
Set-TextValue "D45" "61.29"
Set-TextValue "E45" "  -4.39%  "
Set-TextValue "D46" "8.78"
Set-TextValue "E46" "  +5.64%  "
Set-TextValue "E47" "  +4.48%  "
Set-TextValue "D48" "99.92"
Set-TextValue "E48" "  -0.73%  "
Set-TextValue "D49" "1.19"
Set-TextValue "E49" "  +0.43%  "
Set-TextValue "D50" "2.503.32"
Set-TextValue "E50" "  +2.88%  "
Set-TextValue "D51" "0.427"
Set-TextValue "E51" "  -4.33%  "
